# Adjusting NOV13/14 observation locations.
#
# Rows 11 and 12 on Sheet1 hold an (agency, station-id) pair each.
# The USACE/01440 pair (previously row 11) and the USGS/<station> pair
# (previously row 12) swap places, and the USGS station id text itself
# is corrected from "291929089562600" to "07380260".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: USGS station that used to live in row 12, with its
# corrected station id.
$ws.Range("A11").Value = "USGS"
$ws.Range("B11").Value = "07380260"

# New row 12: USACE station that used to live in row 11.
$ws.Range("A12").Value = "USACE"
$ws.Range("B12").Value = "01440"

# Reflect the edited range as the active selection, like Excel would
# after editing these cells.
$ws.Range("A11:B12").Select()
